$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of the existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I and J columns with data for rows 2-35
$values = @(
    @(6,6),
    @(6,6),
    @(8,8),
    @(5,5),
    @(9,9),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,9),
    @(8,8),
    @(5,6),
    @(7,7),
    @(5,5),
    @(7,7),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,7),
    @(6,6),
    @(8,8),
    @(5,6),
    @(6,6),
    @(6,6),
    @(3,3),
    @(4,4),
    @(5,5),
    @(4,4),
    @(4,4),
    @(4,4)
)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
